$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.883.78"
$ws.Range("E2").Value = "  -0.66%  "
$ws.Range("D3").Value = "'1.907.84"
$ws.Range("E3").Value = "  -0.16%  "
$ws.Range("D4").Value = "'1.004"
$ws.Range("E4").Value = "  -0.11%  "
$ws.Range("D5").Value = "'313.47"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("E6").Value = "  -0.26%  "
$ws.Range("D7").Value = "'0.5002"
$ws.Range("E7").Value = "  +3.80%  "
$ws.Range("D8").Value = "'0.3818"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "'0.07296"
$ws.Range("D10").Value = "'0.9083"
$ws.Range("E10").Value = "  -2.88%  "
$ws.Range("D11").Value = "'20.90"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'0.07674"
$ws.Range("D13").Value = "'1.917.86"
$ws.Range("E13").Value = "  +0.43%  "
$ws.Range("E14").Value = "  -0.45%  "
$ws.Range("D15").Value = "'91.71"
$ws.Range("E15").Value = "  -0.09%  "
$ws.Range("D16").Value = "'1.004"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "'0.000008724"
$ws.Range("E17").Value = "  -1.25%  "
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "  -0.06%  "
$ws.Range("D19").Value = "'27.918.30"
$ws.Range("E19").Value = "  -0.65%  "
$ws.Range("D20").Value = "'14.59"
$ws.Range("E20").Value = "  -1.81%  "
$ws.Range("D21").Value = "'5.181"
$ws.Range("E21").Value = "  +0.06%  "
$ws.Range("E22").Value = "  -0.89%  "
$ws.Range("D23").Value = "'6.600"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Value = "'154.13"
$ws.Range("E24").Value = "  -1.27%  "
$ws.Range("D25").Value = "'1.881"
$ws.Range("E25").Value = "  -2.14%  "
$ws.Range("E26").Value = "  +5.55%  "
$ws.Range("D27").Value = "'18.41"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("D28").Value = "'115.43"
$ws.Range("E28").Value = "  -1.02%  "
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "'0.08990"
$ws.Range("D31").Value = "'3.203"
$ws.Range("E31").Value = "  -3.05%  "
$ws.Range("D32").Value = "'1.235"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").Value = "'0.7658"
$ws.Range("E33").Value = "  -1.49%  "
$ws.Range("D34").Value = "'4.661"
$ws.Range("E34").Value = "  -0.45%  "
$ws.Range("D35").Value = "'0.02063"
$ws.Range("E35").Value = "  +0.25%  "
$ws.Range("D36").Value = "'2.549"
$ws.Range("E36").Value = "  -4.05%  "
$ws.Range("D39").Value = "'3.023"
$ws.Range("E39").Value = "  +1.28%  "
$ws.Range("D40").Value = "'0.05259"
$ws.Range("E40").Value = "  -1.02%  "
$ws.Range("E41").Value = "  -0.79%  "
$ws.Range("D42").Value = "'8.514"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("E43").Value = "  -1.05%  "
$ws.Range("D44").Value = "'111.49"
$ws.Range("E44").Value = "  +2.77%  "
$ws.Range("D45").Value = "'10.61"
$ws.Range("E45").Value = "  -0.52%  "
$ws.Range("D46").Value = "'0.4828"
$ws.Range("E46").Value = "  -0.34%  "
$ws.Range("D47").Value = "'1.002"
$ws.Range("E47").Value = "  -0.36%  "
$ws.Range("D48").Value = "'1.630"
$ws.Range("E48").Value = "  -1.47%  "
$ws.Range("D49").Value = "'67.61"
$ws.Range("E49").Value = "  -0.73%  "
$ws.Range("D50").Value = "'0.06067"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'0.9024"
$ws.Range("E51").Value = "  +0.31%  "

# Row 37/38 swap (TrustWalletToken <-> TheSandbox) with updated values
$ws.Range("B37").Value = "TheSandbox"
$ws.Range("C37").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D37").Value = "'0.5580"
$ws.Range("E37").Value = "  +1.85%  "

$ws.Range("B38").Value = "TrustWalletToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D38").Value = "'1.093"
$ws.Range("E38").Value = "  -1.59%  "
